# "Got Specs from Desktop"
#
# - Record new Desktop PC benchmark timings (F5:F7) on "PythonA1 time
#   measures" as numbers (thousands-separator format) instead of the old
#   placeholder text ("1,763" / "7,6" / "30,206" -> OoT-style strings).
# - Add a new "Desktop Specs" sheet (between "PythonA1 time measures" and
#   "Laptop Specs") with the CPU / Memory / Graphics Card of the desktop
#   machine the new measurements were taken on.

$wb = $excel.ActiveWorkbook

# --- Update the Desktop PC timing column on the measures sheet ---------
$wsMeasures = $wb.Worksheets.Item(1)

$wsMeasures.Range("F5").Value = 1956
$wsMeasures.Range("F5").NumberFormat = "#,##0"

$wsMeasures.Range("F6").Value = 7825
$wsMeasures.Range("F6").NumberFormat = "#,##0"

$wsMeasures.Range("F7").Value = 31984
$wsMeasures.Range("F7").NumberFormat = "#,##0"

# --- Insert a new "Desktop Specs" sheet right before "Laptop Specs" ----
$wsLaptopBefore = $wb.Worksheets.Item("Laptop Specs")
$wsDesktop = $wb.Worksheets.Add($wsLaptopBefore)
$wsDesktop.Name = "Desktop Specs"

$wsDesktop.Range("C6").Value = "CPU"
$wsDesktop.Range("D6").Value = "AMD Ryzen 7 3700X 8-Core Processor 3600Mhz"

$wsDesktop.Range("C7").Value = "Memory (Gb)"
$wsDesktop.Range("D7").Value = 16

$wsDesktop.Range("C8").Value = "Graphics Card"
$wsDesktop.Range("D8").Value = "NVIDIA GeForce RTX 3060"

# --- Restore per-sheet selection state ----------------------------------
# Re-resolve "Laptop Specs" by name: inserting a sheet shifts positional
# handles, so grab a fresh reference before touching it.
$wsLaptop = $wb.Worksheets.Item("Laptop Specs")

[void]$wsDesktop.Range("D12").Select()

[void]$wsLaptop.Activate()
[void]$wsLaptop.Range("D6").Select()

[void]$wsMeasures.Activate()
[void]$wsMeasures.Range("D19").Select()
